$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 507
$ws.Range("J17").Value = 507
$ws.Range("L17").Value = 1521
$ws.Range("N17").Value = -1857

$ws.Range("H100").Value = 40002596
$ws.Range("I100").Value = 1050
$ws.Range("J100").Value = 66670292
$ws.Range("K100").Value = 1050
$ws.Range("L100").Value = 66670292
$ws.Range("M100").Value = -509
$ws.Range("N100").Value = -66671374

$ws.Range("H112").Value = 993284.6
$ws.Range("J112").Value = 1097756.6
$ws.Range("L112").Value = 3293269.8
$ws.Range("N112").Value = -3295485.8

$ws.Range("H115").Value = 8524.556
$ws.Range("I115").Value = 772.4
$ws.Range("K115").Value = 2317.2
$ws.Range("M115").Value = -750.1999999999998

$ws.Range("H116").Value = 13024612
$ws.Range("I116").Value = 5557761
$ws.Range("J116").Value = 27958312
$ws.Range("K116").Value = 5557761
$ws.Range("L116").Value = 27958312
$ws.Range("M116").Value = -5554319
$ws.Range("N116").Value = -27965196

$ws.Range("H125").Value = 1165296
$ws.Range("I125").Value = 7377.3335
$ws.Range("J125").Value = 1358282.5
$ws.Range("K125").Value = 66396.0015
$ws.Range("L125").Value = 12224542.5
$ws.Range("M125").Value = -63936.0015
$ws.Range("N125").Value = -12229462.5

$ws.Range("H127").Value = 2027.9166
$ws.Range("J127").Value = 3200
$ws.Range("L127").Value = 9600
$ws.Range("N127").Value = -19520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1250
$ws.Range("I97").Value = 1250
$ws.Range("K97").Value = 1250
$ws.Range("M97").Value = -754

$ws.Range("H102").Value = 2588.1875
$ws.Range("I102").Value = 2564.2856
$ws.Range("K102").Value = 2564.2856
$ws.Range("M102").Value = -942.2856000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2308.4285
$ws.Range("I94").Value = 1647.5385
$ws.Range("K94").Value = 1647.5385
$ws.Range("M94").Value = -1196.5385

$ws.Range("H96").Value = 15718.571
$ws.Range("I96").Value = 14006
$ws.Range("K96").Value = 14006
$ws.Range("M96").Value = -11260

$ws.Range("H97").Value = 8970.6
$ws.Range("I97").Value = 6213.25
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 6213.25
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -5222.25
$ws.Range("N97").Value = -21982

$ws.Range("H99").Value = 1637.625
$ws.Range("I99").Value = 1400.3
$ws.Range("J99").Value = 2033.1666
$ws.Range("K99").Value = 1400.3
$ws.Range("L99").Value = 2033.1666
$ws.Range("M99").Value = 97.70000000000005
$ws.Range("N99").Value = -5029.1666

$ws.Range("H100").Value = 18908.25
$ws.Range("J100").Value = 18908.25
$ws.Range("L100").Value = 18908.25
$ws.Range("N100").Value = -21072.25

$ws.Range("H102").Value = 11670
$ws.Range("I102").Value = 11670
$ws.Range("K102").Value = 11670
$ws.Range("M102").Value = -8425

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H104").Value = 34888
$ws.Range("J104").Value = 34888
$ws.Range("L104").Value = 34888
$ws.Range("N104").Value = -41876

$ws.Range("H105").Value = 1359.174
$ws.Range("I105").Value = 1312.4286
$ws.Range("J105").Value = 1850
$ws.Range("K105").Value = 1312.4286
$ws.Range("L105").Value = 1850
$ws.Range("M105").Value = 434.5714
$ws.Range("N105").Value = -5344

$ws.Range("H106").Value = 17400
$ws.Range("J106").Value = 17400
$ws.Range("L106").Value = 17400
$ws.Range("N106").Value = -19924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2236992.5
$ws.Range("I31").Value = 2128.4707
$ws.Range("J31").Value = 5690873.5
$ws.Range("K31").Value = 2128.4707
$ws.Range("L31").Value = 5690873.5
$ws.Range("M31").Value = -1833.4707
$ws.Range("N31").Value = -5691463.5

$ws.Range("H34").Value = 2236992.5
$ws.Range("I34").Value = 2128.4707
$ws.Range("J34").Value = 5690873.5
$ws.Range("K34").Value = 2128.4707
$ws.Range("L34").Value = 5690873.5
$ws.Range("M34").Value = -1926.4707
$ws.Range("N34").Value = -5691277.5

$ws.Range("H99").Value = 10025.6
$ws.Range("I99").Value = 6549.3335
$ws.Range("J99").Value = 15240
$ws.Range("K99").Value = 6549.3335
$ws.Range("L99").Value = 15240
$ws.Range("M99").Value = -5051.3335
$ws.Range("N99").Value = -18236

$ws.Range("H126").Value = 10025.6
$ws.Range("I126").Value = 6549.3335
$ws.Range("J126").Value = 15240
$ws.Range("K126").Value = 19648.0005
$ws.Range("L126").Value = 45720
$ws.Range("M126").Value = -17178.0005
$ws.Range("N126").Value = -50660

$ws.Range("H132").Value = 1154.409
$ws.Range("I132").Value = 764.64105
$ws.Range("J132").Value = 4194.6
$ws.Range("K132").Value = 2293.92315
$ws.Range("L132").Value = 12583.8
$ws.Range("M132").Value = 236.0768500000004
$ws.Range("N132").Value = -17643.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5577640
$ws.Range("I121").Value = 1495.7142
$ws.Range("J121").Value = 11153785
$ws.Range("K121").Value = 4487.142599999999
$ws.Range("L121").Value = 33461355
$ws.Range("M121").Value = -3177.142599999999
$ws.Range("N121").Value = -33463975

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 16667880
$ws.Range("I97").Value = 910.55554
$ws.Range("J97").Value = 41668332
$ws.Range("K97").Value = 910.55554
$ws.Range("L97").Value = 41668332
$ws.Range("M97").Value = -414.55554
$ws.Range("N97").Value = -41669324

$ws.Range("H107").Value = 225.25
$ws.Range("I107").Value = 125.5
$ws.Range("J107").Value = 325
$ws.Range("K107").Value = 125.5
$ws.Range("L107").Value = 325
$ws.Range("M107").Value = 1794.5
$ws.Range("N107").Value = -4165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1722.1
$ws.Range("I40").Value = 1671.7778
$ws.Range("J40").Value = 2175
$ws.Range("K40").Value = 1671.7778
$ws.Range("L40").Value = 2175
$ws.Range("M40").Value = -1535.7778
$ws.Range("N40").Value = -2447

$ws.Range("H107").Value = 1150
$ws.Range("I107").Value = 1150
$ws.Range("K107").Value = 1150
$ws.Range("M107").Value = 770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1557.8276
$ws.Range("I96").Value = 1340.8948
$ws.Range("J96").Value = 1970
$ws.Range("K96").Value = 1340.8948
$ws.Range("M96").Value = 32.10519999999997
$ws.Range("N96").Value = -4716
